$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first 16 data rows (rows 2-17), which correspond to the
# earliest 16 quarterly observations (1984-07-01 through 1988-04-01).
# Excel will shift all remaining rows up, so the data that used to live
# in row 18 becomes the new row 2, and the old row 164 becomes row 148.
$ws.Range("A2:B17").EntireRow.Delete()
